# Weekly price update: insert a new weekly price record as row 205
# (shifting all subsequent records down by one row) on the single
# "Sheet1" worksheet of the "Femacal de La Calera - Zanahoria" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 205; everything currently at
# row 205 onward (through 298) shifts down to 206..299.
$ws.Rows("205:205").Insert()

# Populate the newly inserted row 205 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R mirror the neighboring records for
# this market/product (Femacal de La Calera, Coquimbo, Zanahoria,
# Primera, Chillan, $/saco 20 kilos, 20 kg, Hortaliza).
$ws.Range("A205").Value = 3
$ws.Range("B205").Value = "Femacal de La Calera"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = 44609
$ws.Range("E205").Value = 5
$ws.Range("F205").Value = 100114013
$ws.Range("G205").Value = "Zanahoria"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 140
$ws.Range("K205").Value = 9500
$ws.Range("L205").Value = 10000
$ws.Range("M205").Value = 9786
$ws.Range("N205").Value = "$/saco 20 kilos"
$ws.Range("O205").Value = "Chillán"
$ws.Range("P205").Value = 489
$ws.Range("Q205").Value = 20
$ws.Range("R205").Value = "Hortaliza"
